$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 6493799.5
$ws.Range("J17").Value = 7576066
$ws.Range("L17").Value = 22728198
$ws.Range("N17").Value = -22728534

# Row 40
$ws.Range("H40").Value = 2555.25
$ws.Range("J40").Value = 2666.7856
$ws.Range("L40").Value = 2666.7856
$ws.Range("N40").Value = -3016.7856

# Row 62
$ws.Range("H62").Value = 39434
$ws.Range("I62").Value = 19144
$ws.Range("K62").Value = 19144
$ws.Range("M62").Value = -18520

# Row 65
$ws.Range("H65").Value = 39434
$ws.Range("I65").Value = 19144
$ws.Range("K65").Value = 95720
$ws.Range("M65").Value = -92600

# Row 86
$ws.Range("H86").Value = 47039.9
$ws.Range("I86").Value = 6133.3335
$ws.Range("J86").Value = 64571.285
$ws.Range("K86").Value = 6133.3335
$ws.Range("L86").Value = 64571.285
$ws.Range("M86").Value = -5010.3335
$ws.Range("N86").Value = -66817.285

# Row 89
$ws.Range("H89").Value = 47039.9
$ws.Range("I89").Value = 6133.3335
$ws.Range("J89").Value = 64571.285
$ws.Range("K89").Value = 30666.6675
$ws.Range("L89").Value = 322856.425
$ws.Range("M89").Value = -25050.6675
$ws.Range("N89").Value = -334088.425

# Row 96
$ws.Range("H96").Value = 907.3200000000001
$ws.Range("I96").Value = 741
$ws.Range("J96").Value = 1203
$ws.Range("K96").Value = 2223
$ws.Range("L96").Value = 3609
$ws.Range("M96").Value = -850
$ws.Range("N96").Value = -6355

# Row 97
$ws.Range("H97").Value = 2001.7273
$ws.Range("J97").Value = 1702
$ws.Range("L97").Value = 5106
$ws.Range("N97").Value = -6098

# Row 111
$ws.Range("H111").Value = 1356.8
$ws.Range("I111").Value = 899
$ws.Range("J111").Value = 1471.25
$ws.Range("K111").Value = 2697
$ws.Range("L111").Value = 4413.75
$ws.Range("M111").Value = 370
$ws.Range("N111").Value = -10547.75

# Row 113
$ws.Range("H113").Value = 56000
$ws.Range("I113").Value = 56000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 56000
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -52746

# Row 116
$ws.Range("H116").Value = 1766511.5
$ws.Range("J116").Value = 3745.5
$ws.Range("L116").Value = 3745.5
$ws.Range("N116").Value = -10629.5

# Row 129
$ws.Range("H129").Value = 66668410
$ws.Range("I129").Value = 1492.625
$ws.Range("J129").Value = 142859170
$ws.Range("K129").Value = 4477.875
$ws.Range("L129").Value = 428577510
$ws.Range("M129").Value = 522.125
$ws.Range("N129").Value = -428587510

# Row 132
$ws.Range("H132").Value = 4753.9565
$ws.Range("I132").Value = 5728.6875
$ws.Range("J132").Value = 2526
$ws.Range("K132").Value = 17186.0625
$ws.Range("L132").Value = 7578
$ws.Range("M132").Value = -14656.0625
$ws.Range("N132").Value = -12638

# Row 138
$ws.Range("H138").Value = 8898.161
$ws.Range("I138").Value = 7357.0713
$ws.Range("J138").Value = 10167.294
$ws.Range("K138").Value = 22071.2139
$ws.Range("L138").Value = 30501.882
$ws.Range("M138").Value = -16931.2139
$ws.Range("N138").Value = -40781.882


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 5
$ws.Range("H5").Value = 1583.4286
$ws.Range("I5").Value = 363.66666
$ws.Range("J5").Value = 2498.25
$ws.Range("K5").Value = 363.66666
$ws.Range("L5").Value = 2498.25
$ws.Range("M5").Value = -251.66666
$ws.Range("N5").Value = -2722.25

# Row 43
$ws.Range("H43").Value = 11750.454
$ws.Range("J43").Value = 13708.143
$ws.Range("L43").Value = 13708.143
$ws.Range("N43").Value = -14334.143

# Row 44
$ws.Range("H44").Value = 84986.5
$ws.Range("J44").Value = 84986.5
$ws.Range("L44").Value = 84986.5
$ws.Range("N44").Value = -85962.5

# Row 45
$ws.Range("H45").Value = 14639.308
$ws.Range("I45").Value = 21012.25
$ws.Range("J45").Value = 4442.6
$ws.Range("K45").Value = 21012.25
$ws.Range("L45").Value = 4442.6
$ws.Range("M45").Value = -20635.25
$ws.Range("N45").Value = -5196.6

# Row 46
$ws.Range("H46").Value = 11688.167
$ws.Range("J46").Value = 11028.8
$ws.Range("L46").Value = 11028.8
$ws.Range("N46").Value = -11666.8

# Row 61
$ws.Range("H61").Value = 9991.5
$ws.Range("I61").Value = 12622.546
$ws.Range("K61").Value = 12622.546
$ws.Range("M61").Value = -12410.546

# Row 110
$ws.Range("H110").Value = 2246.842
$ws.Range("I110").Value = 1169.8
$ws.Range("K110").Value = 1169.8
$ws.Range("M110").Value = 875.2

# Row 122
$ws.Range("H122").Value = 1082431.2
$ws.Range("I122").Value = 4879.5713
$ws.Range("J122").Value = 2339574.8
$ws.Range("K122").Value = 14638.7139
$ws.Range("L122").Value = 7018724.399999999
$ws.Range("M122").Value = -12188.7139
$ws.Range("N122").Value = -7023624.399999999

# Row 132
$ws.Range("H132").Value = 3733.258
$ws.Range("I132").Value = 3151.5908
$ws.Range("J132").Value = 5155.1113
$ws.Range("K132").Value = 9454.7724
$ws.Range("L132").Value = 15465.3339
$ws.Range("M132").Value = -6924.7724
$ws.Range("N132").Value = -20525.3339

# Row 136
$ws.Range("H136").Value = 9991.5
$ws.Range("I136").Value = 12622.546
$ws.Range("K136").Value = 37867.638
$ws.Range("M136").Value = -35317.638


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 4
$ws.Range("H4").Value = 1583.4286
$ws.Range("I4").Value = 363.66666
$ws.Range("J4").Value = 2498.25
$ws.Range("K4").Value = 363.66666
$ws.Range("L4").Value = 2498.25
$ws.Range("M4").Value = -248.66666
$ws.Range("N4").Value = -2728.25

# Row 94
$ws.Range("H94").Value = 8367.391
$ws.Range("I94").Value = 10466.138
$ws.Range("K94").Value = 10466.138
$ws.Range("M94").Value = -10015.138

# Row 96
$ws.Range("H96").Value = 17802.334
$ws.Range("I96").Value = 14317.429
$ws.Range("K96").Value = 14317.429
$ws.Range("M96").Value = -11571.429

# Row 99
$ws.Range("H99").Value = 11886.525
$ws.Range("I99").Value = 13378.366
$ws.Range("K99").Value = 13378.366
$ws.Range("M99").Value = -11880.366

# Row 105
$ws.Range("H105").Value = 8814.65
$ws.Range("I105").Value = 11191.333
$ws.Range("J105").Value = 5249.625
$ws.Range("K105").Value = 11191.333
$ws.Range("L105").Value = 5249.625
$ws.Range("M105").Value = -9444.333000000001
$ws.Range("N105").Value = -8743.625

# Row 107
$ws.Range("H107").Value = 2669.4062
$ws.Range("I107").Value = 2741
$ws.Range("J107").Value = 2282.8
$ws.Range("K107").Value = 2741
$ws.Range("L107").Value = 2282.8
$ws.Range("M107").Value = -821
$ws.Range("N107").Value = -6122.8

# Row 134
$ws.Range("H134").Value = 7646.125
$ws.Range("I134").Value = 8277.385
$ws.Range("J134").Value = 4910.6665
$ws.Range("K134").Value = 24832.155
$ws.Range("L134").Value = 14731.9995
$ws.Range("M134").Value = -22297.155
$ws.Range("N134").Value = -19801.9995


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2405.4807
$ws.Range("I31").Value = 1688.3846
$ws.Range("K31").Value = 1688.3846
$ws.Range("M31").Value = -1393.3846

# Row 34
$ws.Range("H34").Value = 2405.4807
$ws.Range("I34").Value = 1688.3846
$ws.Range("K34").Value = 1688.3846
$ws.Range("M34").Value = -1486.3846

# Row 62
$ws.Range("H62").Value = 5449.1904
$ws.Range("J62").Value = 7825.5
$ws.Range("L62").Value = 7825.5
$ws.Range("N62").Value = -9073.5

# Row 65
$ws.Range("H65").Value = 5449.1904
$ws.Range("J65").Value = 7825.5
$ws.Range("L65").Value = 39127.5
$ws.Range("N65").Value = -45367.5

# Row 86
$ws.Range("H86").Value = 20488.9
$ws.Range("I86").Value = 17975
$ws.Range("K86").Value = 17975
$ws.Range("M86").Value = -16852

# Row 89
$ws.Range("H89").Value = 20488.9
$ws.Range("I89").Value = 17975
$ws.Range("K89").Value = 89875
$ws.Range("M89").Value = -84259

# Row 94
$ws.Range("H94").Value = 2417.1177
$ws.Range("I94").Value = 3994.4285
$ws.Range("J94").Value = 1313
$ws.Range("K94").Value = 3994.4285
$ws.Range("L94").Value = 1313
$ws.Range("M94").Value = -3543.4285
$ws.Range("N94").Value = -2215

# Row 99
$ws.Range("H99").Value = 5000000
$ws.Range("I99").Value = 5000000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5000000
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4998502

# Row 126
$ws.Range("H126").Value = 5000000
$ws.Range("I126").Value = 5000000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -14997530

# Row 132
$ws.Range("H132").Value = 14517128
$ws.Range("I132").Value = 17551944
$ws.Range("J132").Value = 101749.75
$ws.Range("K132").Value = 52655832
$ws.Range("L132").Value = 305249.25
$ws.Range("M132").Value = -52653302
$ws.Range("N132").Value = -310309.25


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 17
$ws.Range("H17").Value = 2239.6
$ws.Range("I17").Value = 1302
$ws.Range("J17").Value = 5990
$ws.Range("K17").Value = 3906
$ws.Range("L17").Value = 17970
$ws.Range("M17").Value = -3737
$ws.Range("N17").Value = -18308

# Row 39
$ws.Range("H39").Value = 1761
$ws.Range("J39").Value = 6775
$ws.Range("L39").Value = 20325
$ws.Range("N39").Value = -20913

# Row 122
$ws.Range("H122").Value = 2532.8276
$ws.Range("J122").Value = 2659.577
$ws.Range("L122").Value = 23936.193
$ws.Range("N122").Value = -28836.193

# Row 132
$ws.Range("H132").Value = 27863328
$ws.Range("I132").Value = 1649.5
$ws.Range("J132").Value = 41794170
$ws.Range("K132").Value = 14845.5
$ws.Range("L132").Value = 376147530
$ws.Range("M132").Value = -12315.5
$ws.Range("N132").Value = -376152590


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 9116.929
$ws.Range("I70").Value = 8819.75
$ws.Range("K70").Value = 8819.75
$ws.Range("M70").Value = -8549.75

# Row 73
$ws.Range("H73").Value = 9116.929
$ws.Range("I73").Value = 8819.75
$ws.Range("K73").Value = 8819.75
$ws.Range("M73").Value = -7883.75

# Row 80
$ws.Range("H80").Value = 4835
$ws.Range("I80").Value = 4789
$ws.Range("J80").Value = 4865.6665
$ws.Range("K80").Value = 4789
$ws.Range("L80").Value = 4865.6665
$ws.Range("M80").Value = -3791
$ws.Range("N80").Value = -6861.6665

# Row 83
$ws.Range("H83").Value = 4835
$ws.Range("I83").Value = 4789
$ws.Range("J83").Value = 4865.6665
$ws.Range("K83").Value = 23945
$ws.Range("L83").Value = 24328.3325
$ws.Range("M83").Value = -18953
$ws.Range("N83").Value = -34312.3325

# Row 97
$ws.Range("H97").Value = 8576.037
$ws.Range("I97").Value = 10200.429
$ws.Range("J97").Value = 2890.6667
$ws.Range("K97").Value = 10200.429
$ws.Range("L97").Value = 2890.6667
$ws.Range("M97").Value = -9704.429
$ws.Range("N97").Value = -3882.6667

# Row 107
$ws.Range("H107").Value = 635
$ws.Range("I107").Value = 494
$ws.Range("K107").Value = 494
$ws.Range("M107").Value = 1426

# Row 113
$ws.Range("H113").Value = 3184.0833
$ws.Range("I113").Value = 3111
$ws.Range("J113").Value = 3198.7
$ws.Range("K113").Value = 3111
$ws.Range("L113").Value = 3198.7
$ws.Range("M113").Value = -941
$ws.Range("N113").Value = -7538.7

# Row 122
$ws.Range("H122").Value = 30114.428
$ws.Range("I122").Value = 38854.6
$ws.Range("J122").Value = 8264
$ws.Range("K122").Value = 116563.8
$ws.Range("L122").Value = 24792
$ws.Range("M122").Value = -114113.8
$ws.Range("N122").Value = -29692

# Row 132
$ws.Range("H132").Value = 4416.4546
$ws.Range("I132").Value = 3475.7742
$ws.Range("K132").Value = 10427.3226
$ws.Range("M132").Value = -7897.3226

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 4870.347
$ws.Range("I22").Value = 6301.8076
$ws.Range("J22").Value = 3252.1738
$ws.Range("K22").Value = 6301.8076
$ws.Range("L22").Value = 3252.1738
$ws.Range("M22").Value = -6006.8076
$ws.Range("N22").Value = -3842.1738

# Row 27
$ws.Range("H27").Value = 4870.347
$ws.Range("I27").Value = 6301.8076
$ws.Range("J27").Value = 3252.1738
$ws.Range("K27").Value = 6301.8076
$ws.Range("L27").Value = 3252.1738
$ws.Range("M27").Value = -6194.8076
$ws.Range("N27").Value = -3466.1738

# Row 46
$ws.Range("H46").Value = 2346.7646
$ws.Range("I46").Value = 1530.7
$ws.Range("J46").Value = 3512.5715
$ws.Range("K46").Value = 1530.7
$ws.Range("L46").Value = 3512.5715
$ws.Range("M46").Value = -1342.7
$ws.Range("N46").Value = -3888.5715

# Row 61
$ws.Range("H61").Value = 1711.579
$ws.Range("I61").Value = 1526.3636
$ws.Range("J61").Value = 1966.25
$ws.Range("K61").Value = 1526.3636
$ws.Range("L61").Value = 1966.25
$ws.Range("M61").Value = -1324.3636
$ws.Range("N61").Value = -2370.25

# Row 68
$ws.Range("H68").Value = 3948.4
$ws.Range("I68").Value = 2445
$ws.Range("K68").Value = 2445
$ws.Range("M68").Value = -1696

# Row 71
$ws.Range("H71").Value = 3948.4
$ws.Range("I71").Value = 2445
$ws.Range("K71").Value = 12225
$ws.Range("M71").Value = -8481

# Row 93
$ws.Range("H93").Value = 2221.4092
$ws.Range("I93").Value = 2634.3845
$ws.Range("J93").Value = 1624.8889
$ws.Range("K93").Value = 2634.3845
$ws.Range("L93").Value = 1624.8889
$ws.Range("M93").Value = -1386.3845
$ws.Range("N93").Value = -4120.8889

# Row 113
$ws.Range("H113").Value = 1711.579
$ws.Range("I113").Value = 1526.3636
$ws.Range("J113").Value = 1966.25
$ws.Range("K113").Value = 1526.3636
$ws.Range("L113").Value = 1966.25
$ws.Range("M113").Value = 643.6364000000001
$ws.Range("N113").Value = -6306.25

# Row 122
$ws.Range("H122").Value = 5483.25
$ws.Range("I122").Value = 3799.5
$ws.Range("J122").Value = 5820
$ws.Range("K122").Value = 11398.5
$ws.Range("L122").Value = 17460
$ws.Range("M122").Value = -8948.5
$ws.Range("N122").Value = -22360

# Row 132
$ws.Range("H132").Value = 19454.53
$ws.Range("I132").Value = 21748.23
$ws.Range("K132").Value = 65244.69
$ws.Range("M132").Value = -62714.69

# Row 136
$ws.Range("H136").Value = 3910.25
$ws.Range("I136").Value = 2057.5557
$ws.Range("J136").Value = 4787.8423
$ws.Range("K136").Value = 6172.6671
$ws.Range("L136").Value = 14363.5269
$ws.Range("M136").Value = -3622.6671
$ws.Range("N136").Value = -19463.5269


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 18
$ws.Range("H18").Value = 14249.75
$ws.Range("I18").Value = 14999
$ws.Range("J18").Value = 14000
$ws.Range("K18").Value = 14999
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = -14826
$ws.Range("N18").Value = -14346

# Row 96
$ws.Range("H96").Value = 22808198
$ws.Range("I96").Value = 10001353
$ws.Range("J96").Value = 37038024
$ws.Range("K96").Value = 10001353
$ws.Range("L96").Value = 37038024
$ws.Range("M96").Value = -9999980
$ws.Range("N96").Value = -37040770

# Row 113
$ws.Range("H113").Value = 2937.9285
$ws.Range("J113").Value = 6997
$ws.Range("L113").Value = 20991
$ws.Range("N113").Value = -25331

# Row 122
$ws.Range("H122").Value = 8763.058999999999
$ws.Range("I122").Value = 6664.1113
$ws.Range("K122").Value = 19992.3339
$ws.Range("M122").Value = -17542.3339

# Row 132
$ws.Range("H132").Value = 14117.546
$ws.Range("I132").Value = 16295.68
$ws.Range("J132").Value = 7310.875
$ws.Range("K132").Value = 48887.04
$ws.Range("L132").Value = 21932.625
$ws.Range("M132").Value = -46357.04
$ws.Range("N132").Value = -26992.625

# Row 140
$ws.Range("H140").Value = 70103.14
$ws.Range("J140").Value = 70103.14
$ws.Range("L140").Value = 70103.14
$ws.Range("N140").Value = -80463.14

# Row 141
$ws.Range("H141").Value = 111105.5
$ws.Range("J141").Value = 111105.5
$ws.Range("L141").Value = 111105.5
$ws.Range("N141").Value = -121465.5

